$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.436.46'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.571.20'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.43'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.495'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.25%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.05'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0864'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.796.90'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.585.25'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.521'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.12'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.469.84'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '213.69'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.69%  '
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.26'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.12'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.77'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.41%  '
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.89'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.79'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  -0.92%  '
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('E32').Value = '  -1.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.358.60'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.970'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.31'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.531'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.26'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.710.65'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.35'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0995'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0494'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.58%  '

Write-Host "Updated cryptos list values"
